# "UI with Cable Sizing"
# - add a new "Sheet1" tab (after "in") that holds a scratch list of cable
#   sizes that aren't used on the main table (2C#4, 2C#6, 2C#9, 2C#14)
# - rebuild the "in" sheet's table: drop the Radius / LB Per Foot columns,
#   reorder to Size | Diameter | Lb Per Foot | Cross Sect. Area, bold the
#   header row, and repopulate with the full cable-size list, computing
#   Cross Sect. Area with a formula instead of a hard-coded number.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- new "Sheet1" tab, placed right after "in" ------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "2C#4"
$ws2.Range("A2").Value = "2C#6"
$ws2.Range("A3").Value = "2C#9"
$ws2.Range("A4").Value = "2C#14"

# --- rebuild the "in" sheet --------------------------------------------
$ws1.Cells.Clear()

$ws1.Columns.Item(1).ColumnWidth = 12.43
$ws1.Columns.Item(2).ColumnWidth = 9.29
$ws1.Columns.Item(3).ColumnWidth = 10.86
$ws1.Columns.Item(4).ColumnWidth = 15.14

# header row
$ws1.Range("A1").Value = "Size"
$ws1.Range("B1").Value = "Diameter"
$ws1.Range("D1").Value = "Cross Sect. Area"
$ws1.Range("A1:D1").Font.Bold = $true

# data rows: Size, Diameter, Lb Per Foot, Cross Sect. Area (formula)
$data = @(
    @(2,  "STAR QUAD", 0.9,     0.86),
    @(3,  "3C#4",      1.3,     0.96),
    @(4,  "3C#6",      1.175,   0.735),
    @(5,  "7C#14",     0.99,    0.495),
    @(6,  "7C#12",     1.485,   0.743),
    @(7,  "7C#10",     1.485,   0.743),
    @(8,  "10C#14",    1.232,   0.711),
    @(9,  "12C#14",    1.27,    0.767),
    @(10, "14C#14",    1.332,   0.865),
    @(11, "19C#14",    1.507,   1.284),
    @(12, "25C#14",    1.24,    0.88)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws1.Range("A$r").Value = $row[1]
    $ws1.Range("B$r").Value = $row[2]
    $ws1.Range("C$r").Value = $row[3]
}

# Cross Sect. Area formulas - rows 3:7 share one formula group (mirrors
# how Excel groups a single fill/copy down the column), the rest are
# entered individually.
$ws1.Range("D2").Formula = "=ROUND(PI() * ((B2/2)^2), 2)"
$ws1.Range("D3:D7").Formula = "=ROUND(PI() * ((B3/2)^2), 2)"
$ws1.Range("D8").Formula = "=ROUND(PI() * ((B8/2)^2), 2)"
$ws1.Range("D9").Formula = "=ROUND(PI() * ((B9/2)^2), 2)"
$ws1.Range("D10").Formula = "=ROUND(PI() * ((B10/2)^2), 2)"
$ws1.Range("D11").Formula = "=ROUND(PI() * ((B11/2)^2), 2)"
$ws1.Range("D12").Formula = "=ROUND(PI() * ((B12/2)^2), 2)"

# column header text that must land LAST in the shared-string table
$ws1.Range("C1").Value = "Lb Per Foot"

$ws1.PageSetup.Orientation = 1

# selection / active sheet bookkeeping - touch Sheet1 first so that the
# final selection (and active tab) ends up back on "in"
$ws2.Range("A1:A4").Select()
$ws1.Range("E5").Select()
